$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price / volume / coin reorderings)

# Row 2
$ws.Range("D2").Value = "64.213.33"

# Row 3
$ws.Range("D3").Value = "3.174.65"
$ws.Range("E3").Value = "  -8.32%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'564.51"
$ws.Range("E5").Value = "  -3.33%  "

# Row 6
$ws.Range("D6").Value = "'168.98"
$ws.Range("E6").Value = "  -5.08%  "

# Row 7
$ws.Range("D7").Value = "'0.608"
$ws.Range("E7").Value = "  -3.48%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "3.173.03"
$ws.Range("E9").Value = "  -8.34%  "

# Row 10
$ws.Range("E10").Value = "  -7.08%  "

# Row 11
$ws.Range("D11").Value = "'6.61"
$ws.Range("E11").Value = "  -5.13%  "

# Row 12
$ws.Range("E12").Value = "  -5.60%  "

# Row 13
$ws.Range("D13").Value = "3.728.74"
$ws.Range("E13").Value = "  -8.22%  "

# Row 14
$ws.Range("E14").Value = "  +1.35%  "

# Row 15
$ws.Range("D15").Value = "'27.36"
$ws.Range("E15").Value = "  -9.45%  "

# Row 16
$ws.Range("D16").Value = "64.215.47"
$ws.Range("E16").Value = "  -3.14%  "

# Row 17
$ws.Range("D17").Value = "'0.0000163"
$ws.Range("E17").Value = "  -5.42%  "

# Row 18
$ws.Range("D18").Value = "3.175.33"
$ws.Range("E18").Value = "  -8.41%  "

# Row 19
$ws.Range("E19").Value = "  -4.32%  "

# Row 20
$ws.Range("D20").Value = "'12.97"
$ws.Range("E20").Value = "  -6.51%  "

# Row 21
$ws.Range("D21").Value = "'352.64"
$ws.Range("E21").Value = "  -5.11%  "

# Row 22
$ws.Range("E22").Value = "  -6.43%  "

# Row 23
$ws.Range("E23").Value = "  +0.18%  "

# Row 24
$ws.Range("D24").Value = "'68.59"
$ws.Range("E24").Value = "  -6.41%  "

# Row 25
$ws.Range("D25").Value = "'0.504"
$ws.Range("E25").Value = "  -6.15%  "

# Row 26
$ws.Range("D26").Value = "'0.0000118"
$ws.Range("E26").Value = "  -5.82%  "

# Row 27
$ws.Range("D27").Value = "'9.59"
$ws.Range("E27").Value = "  -4.59%  "

# Row 28
$ws.Range("E28").Value = "  -0.74%  "

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30
$ws.Range("E30").Value = "  -0.12%  "

# Row 31
$ws.Range("E31").Value = "  -7.15%  "

# Row 32
$ws.Range("E32").Value = "  -5.06%  "

# Row 33
$ws.Range("D33").Value = "'21.96"
$ws.Range("E33").Value = "  -7.34%  "

# Row 34
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'6.64"
$ws.Range("E34").Value = "  -6.28%  "

# Row 35
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.20"
$ws.Range("E35").Value = "  -5.46%  "

# Row 36
$ws.Range("D36").Value = "'1.43"
$ws.Range("E36").Value = "  -8.16%  "

# Row 37
$ws.Range("D37").Value = "'154.38"
$ws.Range("E37").Value = "  -4.18%  "

# Row 38
$ws.Range("E38").Value = "  -7.76%  "

# Row 39
$ws.Range("D39").Value = "'26.03"
$ws.Range("E39").Value = "  -7.14%  "

# Row 40
$ws.Range("E40").Value = "  -6.82%  "

# Row 41
$ws.Range("D41").Value = "'2.48"
$ws.Range("E41").Value = "  -4.51%  "

# Row 42
$ws.Range("D42").Value = "2.614.95"
$ws.Range("E42").Value = "  -7.13%  "

# Row 43
$ws.Range("E43").Value = "  -7.76%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'5.99"
$ws.Range("E44").Value = "  -7.28%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'39.30"
$ws.Range("E45").Value = "  -1.83%  "

# Row 46
$ws.Range("D46").Value = "'0.0646"
$ws.Range("E46").Value = "  -6.98%  "

# Row 47
$ws.Range("D47").Value = "'23.64"
$ws.Range("E47").Value = "  -6.39%  "

# Row 48
$ws.Range("D48").Value = "'317.32"
$ws.Range("E48").Value = "  -7.42%  "

# Row 49
$ws.Range("E49").Value = "  -7.46%  "

# Row 50
$ws.Range("E50").Value = "  -3.83%  "

# Row 51
$ws.Range("E51").Value = "  -0.05%  "
